$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "S.No."
$ws.Range("B1").Value = "ScenarioName"
$ws.Range("C1").Value = "resource"
$ws.Range("D1").Value = "api_Request"
$ws.Range("E1").Value = "key"
$ws.Range("F1").Value = "content_Type"
$ws.Range("G1").Value = "accuracy"
$ws.Range("H1").Value = "name"
$ws.Range("I1").Value = "phone_number"
$ws.Range("J1").Value = "address"
$ws.Range("K1").Value = "website"
$ws.Range("L1").Value = "language"
$ws.Range("M1").Value = "lat"
$ws.Range("N1").Value = "lng"
$ws.Range("O1").Value = "types"

# --- Row 2 ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Verify if place is being added using Add Place API"
$ws.Range("C2").Value = "/maps/api/place/add/json"
$ws.Range("D2").Value = "POST"
$ws.Range("E2").Value = "qaclick123"
$ws.Range("F2").Value = "JSON"
$ws.Range("G2").Value = "'50"
$ws.Range("H2").Value = "Frontline house"
$ws.Range("I2").Value = "'9822789334"
$ws.Range("J2").Value = "India"
$ws.Range("K2").Value = "pranjal.com"
$ws.Range("L2").Value = "Hindi"
$ws.Range("M2").Value = "'-38.383494"
$ws.Range("N2").Value = "'33.427362"
$ws.Range("O2").Value = "Shop#Shoe Park"

# --- Row 3 ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Verify if place is not being added using Add Place API"
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "qaclick123"
$ws.Range("F3").Value = "JSON"
$ws.Range("G3").Value = "NA"
$ws.Range("H3").Value = "NA"
$ws.Range("I3").Value = "NA"
$ws.Range("J3").Value = "NA"
$ws.Range("K3").Value = "NA"
$ws.Range("L3").Value = "NA"
$ws.Range("M3").Value = "NA"
$ws.Range("N3").Value = "NA"
$ws.Range("O3").Value = "NA"

# --- Rows 4-17: column A only, blank but styled ---
for ($r = 4; $r -le 17; $r++) {
    $ws.Range("A$r").Value = ""
}

# --- Styles: center+middle alignment on column A (rows 2-17) ---
# Build the combined style once on a scratch cell, then copy just the
# formatting over so the engine doesn't retain an intermediate
# "horizontal-only" style in the workbook's style table.
$ws.Range("Z1").HorizontalAlignment = -4108
$ws.Range("Z1").VerticalAlignment = -4108
$ws.Range("Z1").Copy()
$ws.Range("A2:A17").PasteSpecial(-4122)
$ws.Range("Z1").Clear()

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 48.1666666666667
$ws.Columns.Item(3).ColumnWidth = 24.0221354166667
$ws.Columns.Item(4).ColumnWidth = 24.0221354166667
$ws.Columns.Item(5).ColumnWidth = 9.16666666666667
$ws.Columns.Item(6).ColumnWidth = 11.3072916666667
$ws.Columns.Item(9).ColumnWidth = 10.1666666666667

# --- Selection ---
$ws.Range("D3").Select()
